$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header/index row (row 1) for columns B:E
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Update row 2 values (meanEMG / passive torque data, common max ROM)
$ws.Range("B2").Value = 24.108488991975815
$ws.Range("C2").Value = 16.403585443062106
$ws.Range("D2").Value = 32.89271844946753
$ws.Range("E2").Value = 16.859891899922577

# Update row 3 values
$ws.Range("B3").Value = 32.207610945824023
$ws.Range("C3").Value = 19.470153128188006
$ws.Range("D3").Value = 36.412500352069117
$ws.Range("E3").Value = 16.016676821730766

# Restore the original selection on the sheet (B1:E3)
$ws.Range("B1:E3").Select()
